$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change D1 from "#unique=false" to "#match=all"
$ws.Range("D1").Value = "#match=all"

# Update selection to D2, matching the saved selection in the diff
$ws.Range("D2").Select()
